# BP-1588 category all apis
#
# Insert a new "category" column into the "Курсы" (Courses) sheet, between
# the existing "learningOutcome" column (D) and "type" column (old E).
# The new column becomes E ("category" / "Профориентация"); the old E/F
# columns (type/options) shift right to F/G.
#
# Also: the "Курсы" sheet becomes the active tab (it was "О компании"
# before), reflected by selecting a range on it at the end of the script.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Курсы")

# Shift old column E ("type") and F ("options") one slot to the right and
# open up a new column E for "category". Match its width to the existing
# "fullDescription"/"learningOutcome" columns (C:D).
$ws.Columns.Item(5).Insert()
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# Header cell: "category" (same look as the other header cells — left /
# centered text).
$ws.Range("E1").Value = "category"
$ws.Range("E1").HorizontalAlignment = -4131   # xlLeft
$ws.Range("E1").VerticalAlignment = -4108     # xlCenter

# Data cell: "Профориентация" — wrapped, vertically justified text.
$ws.Range("E2").Value = "Профориентация"
$ws.Range("E2").WrapText = $true
$ws.Range("E2").VerticalAlignment = -4130     # xlJustify

# Make "Курсы" the active sheet / selected range (matches the new
# tabSelected + selection in the saved workbook).
$ws.Range("E1:E2").Select() | Out-Null

Write-Output "done"
